# Update Excel Download with SSZ Logo
#
# The title sheet ("Tabelle1") used to carry a standalone text row reading
# "Statistik Stadt Zürich" right under the header row. Now that the export
# shows the SSZ logo (image) instead, that text row is redundant, so it is
# removed outright - the rows beneath it (address block, "Erstellt am" /
# "Datum", "Inhalt", "T_1" placeholder) all shift up by one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Delete row 2 ("Statistik Stadt Zürich") entirely; everything below moves
# up by one row (old row 3 -> new row 2, old row 12 -> new row 11, etc.).
$ws.Rows.Item(2).Delete()

# Leave the selection on the row that is now in position 2 (matches the
# workbook's saved selection state after the edit).
$ws.Rows.Item(2).Select()

# Refresh the recorded project folder (Microsoft/x15ac absolute-path hint)
# now that the workbook lives in the new "LIMA_2.0" folder rather than the
# old "LIMAneu" one.
$wb.AbsPath = "O:\Projekte\Bodenpreise\6_Shiny\LIMA_2.0\"
